# "Generate Report for Handoff"
# Regenerating the localization-status report refreshes the "Latest Handoff
# Datetime" value for the 85ba0f79-... source file in both the zh-cn and
# de-de target-language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-08 12:27:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-08 12:27:25"
